$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the image path for the PESEL document card as a new cell (P2),
# matching the "documentType"/"image" column header in row 1.
$ws.Range("P2").Value = "C:Users/vano/Documents/GitHub/ZPI_VAF/iaff_assistant/images/Cards/PESEL.jpg"

# Leave the selection on the newly-populated cell, as happens naturally
# after typing a value into it in the Excel UI.
[void]$ws.Range("P2").Select()
